$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D (Price) and E (Volume 1h) columns for rows 2-50 to match latest scrape
$ws.Range("D2").Value2 = "'27.460.74"
$ws.Range("E2").Value2 = "  +4.96%  "
$ws.Range("D3").Value2 = "'1.815.43"
$ws.Range("E3").Value2 = "  +5.68%  "
$ws.Range("D4").Value2 = "'1.001"
$ws.Range("D5").Value2 = "'343.02"
$ws.Range("E5").Value2 = "  +2.95%  "
$ws.Range("E6").Value2 = "  -0.03%  "
$ws.Range("D7").Value2 = "'0.3810"
$ws.Range("E7").Value2 = "  +3.09%  "
$ws.Range("D8").Value2 = "'0.3495"
$ws.Range("E8").Value2 = "  +4.52%  "
$ws.Range("D9").Value2 = "'49.02"
$ws.Range("E9").Value2 = "  -0.54%  "
$ws.Range("D10").Value2 = "'1.232"
$ws.Range("E10").Value2 = "  +3.57%  "
$ws.Range("D11").Value2 = "'0.07746"
$ws.Range("E11").Value2 = "  +3.72%  "
$ws.Range("D12").Value2 = "'0.9996"
$ws.Range("E12").Value2 = "  -0.03%  "
$ws.Range("D13").Value2 = "'22.21"
$ws.Range("E13").Value2 = "  +10.35%  "
$ws.Range("D14").Value2 = "'6.598"
$ws.Range("E14").Value2 = "  +4.16%  "
$ws.Range("D15").Value2 = "'1.815.03"
$ws.Range("E15").Value2 = "  +5.43%  "
$ws.Range("D16").Value2 = "'7.222"
$ws.Range("E16").Value2 = "  +4.03%  "
$ws.Range("E17").Value2 = "  +3.76%  "
$ws.Range("D18").Value2 = "'0.06712"
$ws.Range("E18").Value2 = "  +0.90%  "
$ws.Range("D19").Value2 = "'86.16"
$ws.Range("E19").Value2 = "  +5.11%  "
$ws.Range("D20").Value2 = "'0.9998"
$ws.Range("E20").Value2 = "  +0.08%  "
$ws.Range("D21").Value2 = "'17.62"
$ws.Range("E21").Value2 = "  +7.31%  "
$ws.Range("D22").Value2 = "'6.588"
$ws.Range("E22").Value2 = "  +8.14%  "
$ws.Range("E23").Value2 = "  +1.53%  "
$ws.Range("D24").Value2 = "'27.449.56"
$ws.Range("E24").Value2 = "  +5.17%  "
$ws.Range("D25").Value2 = "'2.467"
$ws.Range("E25").Value2 = "  -0.35%  "
$ws.Range("D26").Value2 = "'2.671"
$ws.Range("E26").Value2 = "  +8.87%  "
$ws.Range("D27").Value2 = "'22.02"
$ws.Range("E27").Value2 = "  +14.33%  "
$ws.Range("D28").Value2 = "'1.466"
$ws.Range("E28").Value2 = "  +6.71%  "
$ws.Range("D29").Value2 = "'153.86"
$ws.Range("E29").Value2 = "  +1.31%  "
$ws.Range("D30").Value2 = "'2.022.30"
$ws.Range("E30").Value2 = "  +5.66%  "
$ws.Range("D31").Value2 = "'135.80"
$ws.Range("E31").Value2 = "  +4.91%  "
$ws.Range("D32").Value2 = "'6.341"
$ws.Range("E32").Value2 = "  +6.34%  "
$ws.Range("E33").Value2 = "  -1.63%  "
$ws.Range("D34").Value2 = "'13.93"
$ws.Range("E34").Value2 = "  +7.84%  "
$ws.Range("D35").Value2 = "'0.08784"
$ws.Range("E35").Value2 = "  +2.58%  "
$ws.Range("D36").Value2 = "'1.689"
$ws.Range("E36").Value2 = "  -1.24%  "
$ws.Range("D37").Value2 = "'5.617"
$ws.Range("E37").Value2 = "  +4.73%  "
$ws.Range("D38").Value2 = "'0.6979"
$ws.Range("E38").Value2 = "  +12.93%  "
$ws.Range("D39").Value2 = "'0.2268"
$ws.Range("E39").Value2 = "  +5.99%  "
$ws.Range("D40").Value2 = "'0.02403"
$ws.Range("E40").Value2 = "  +3.07%  "
$ws.Range("D41").Value2 = "'0.06476"
$ws.Range("E41").Value2 = "  +3.80%  "
$ws.Range("D42").Value2 = "'8.946"
$ws.Range("E42").Value2 = "  +4.14%  "
$ws.Range("D43").Value2 = "'1.294"
$ws.Range("E43").Value2 = "  +4.76%  "
$ws.Range("D44").Value2 = "'14.72"
$ws.Range("E44").Value2 = "  +1.99%  "
$ws.Range("D45").Value2 = "'0.6541"
$ws.Range("E45").Value2 = "  +10.69%  "
$ws.Range("D46").Value2 = "'0.9989"
$ws.Range("E46").Value2 = "  -0.09%  "
$ws.Range("D47").Value2 = "'4.015"
$ws.Range("E47").Value2 = "  +4.67%  "
$ws.Range("D48").Value2 = "'2.178"
$ws.Range("E48").Value2 = "  +7.71%  "
$ws.Range("D49").Value2 = "'133.38"
$ws.Range("E49").Value2 = "  +3.59%  "
$ws.Range("D50").Value2 = "'0.07323"
$ws.Range("E50").Value2 = "  +0.84%  "

# Row 51: coin changed from Aave to Stacks
$ws.Range("B51").Value2 = "Stacks"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value2 = "'1.273"
$ws.Range("E51").Value2 = "  +20.03%  "
